$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '34.226.08'
$ws.Range('E2').Value = '  -1.03%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.783.72'
$ws.Range('E3').Value = '  -2.60%  '
$ws.Range('E4').Value = '  +0.21%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '224.12'
$ws.Range('E5').Value = '  -2.83%  '
$ws.Range('E6').Value = '  +0.33%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '31.98'
$ws.Range('E8').Value = '  +0.92%  '
$ws.Range('E9').Value = '  -1.34%  '
$ws.Range('E10').Value = '  -2.30%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0931'
$ws.Range('E11').Value = '  +0.11%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.039.81'
$ws.Range('E12').Value = '  -2.46%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.17'
$ws.Range('E13').Value = '  +6.91%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.782.71'
$ws.Range('E14').Value = '  -2.32%  '
$ws.Range('E15').Value = '  -3.80%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '34.233.88'
$ws.Range('E16').Value = '  -0.73%  '
$ws.Range('E17').Value = '  -1.55%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '68.76'
$ws.Range('E18').Value = '  -1.36%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '254.09'
$ws.Range('E19').Value = '  -2.32%  '
$ws.Range('E20').Value = '  -1.87%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.999'
$ws.Range('E21').Value = '  +0.49%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.35'
$ws.Range('E22').Value = '  -2.09%  '
$ws.Range('E23').Value = '  -3.75%  '
$ws.Range('E24').Value = '  -3.31%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '157.27'
$ws.Range('E25').Value = '  -0.47%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '16.37'
$ws.Range('E26').Value = '  -2.03%  '
$ws.Range('E27').Value = '  -1.76%  '
$ws.Range('E28').Value = '  -1.45%  '
$ws.Range('E29').Value = '  +0.15%  '
$ws.Range('E30').Value = '  -3.49%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0513'
$ws.Range('E31').Value = '  -1.34%  '
$ws.Range('E32').Value = '  -2.72%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.59'
$ws.Range('E33').Value = '  +0.06%  '
$ws.Range('E34').Value = '  +3.57%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.442.40'
$ws.Range('E35').Value = '  -7.21%  '
$ws.Range('E36').Value = '  -2.77%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0187'
$ws.Range('E37').Value = '  -1.35%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.623'
$ws.Range('E38').Value = '  -2.00%  '
$ws.Range('E39').Value = '  +1.48%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '82.94'
$ws.Range('E40').Value = '  -2.55%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.34'
$ws.Range('E41').Value = '  +0.45%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.889'
$ws.Range('E42').Value = '  -3.11%  '
$ws.Range('E43').Value = '  -5.29%  '
$ws.Range('E44').Value = '  -3.06%  '
$ws.Range('E45').Value = '  -1.94%  '
$ws.Range('E46').Value = '  +0.38%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.940.27'
$ws.Range('E47').Value = '  -2.57%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '12.19'
$ws.Range('E48').Value = '  -2.19%  '
$ws.Range('E49').Value = '  -0.05%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '98.40'
$ws.Range('E50').Value = '  +0.48%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '49.32'
$ws.Range('E51').Value = '  -6.72%  '
